{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer-key strings in the\n// table with their new values (see commit diff). Each old string is\n// unique within the document, so a direct search+replace per pair is\n// safe and order-independent.\nconst replacements = [\n  [\"102\u00f73=34, 0\", \"112\u00f78=14, 0\"],\n  [\"441\u00f72=220, 1\", \"428\u00f79=47, 5\"],\n  [\"223\u00f74=55, 3\", \"289\u00f78=36, 1\"],\n  [\"592\u00f77=84, 4\", \"262\u00f74=65, 2\"],\n  [\"610\u00f73=203, 1\", \"110\u00f76=18, 2\"],\n  [\"657\u00f76=109, 3\", \"602\u00f77=86, 0\"],\n  [\"333\u00f76=55, 3\", \"706\u00f73=235, 1\"],\n  [\"130\u00f78=16, 2\", \"106\u00f72=53, 0\"],\n  [\"338\u00f77=48, 2\", \"785\u00f74=196, 1\"],\n  [\"646\u00f79=71, 7\", \"995\u00f72=497, 1\"],\n  [\"488\u00f77=69, 5\", \"113\u00f78=14, 1\"],\n  [\"491\u00f73=163, 2\", \"234\u00f78=29, 2\"],\n  [\"232\u00f77=33, 1\", \"429\u00f75=85, 4\"],\n  [\"215\u00f79=23, 8\", \"224\u00f75=44, 4\"],\n  [\"256\u00f74=64, 0\", \"266\u00f77=38, 0\"],\n  [\"376\u00f73=125, 1\", \"851\u00f79=94, 5\"],\n  [\"550\u00f75=110, 0\", \"556\u00f77=79, 3\"],\n  [\"519\u00f74=129, 3\", \"154\u00f74=38, 2\"],\n  [\"235\u00f78=29, 3\", \"181\u00f72=90, 1\"],\n  [\"537\u00f76=89, 3\", \"811\u00f74=202, 3\"],\n  [\"838\u00f75=167, 3\", \"197\u00f76=32, 5\"],\n  [\"723\u00f75=144, 3\", \"226\u00f72=113, 0\"],\n  [\"675\u00f79=75, 0\", \"421\u00f75=84, 1\"],\n  [\"122\u00f78=15, 2\", \"378\u00f75=75, 3\"],\n  [\"535\u00f79=59, 4\", \"118\u00f77=16, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" answer-key strings in the\n# table with their new values (see commit diff). Each old string is\n# unique within the document, so a direct Find/Replace per pair is\n# safe and order-independent.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Old = \"102\u00f73=34, 0\";   New = \"112\u00f78=14, 0\" },\n    @{ Old = \"441\u00f72=220, 1\";  New = \"428\u00f79=47, 5\" },\n    @{ Old = \"223\u00f74=55, 3\";   New = \"289\u00f78=36, 1\" },\n    @{ Old = \"592\u00f77=84, 4\";   New = \"262\u00f74=65, 2\" },\n    @{ Old = \"610\u00f73=203, 1\";  New = \"110\u00f76=18, 2\" },\n    @{ Old = \"657\u00f76=109, 3\";  New = \"602\u00f77=86, 0\" },\n    @{ Old = \"333\u00f76=55, 3\";   New = \"706\u00f73=235, 1\" },\n    @{ Old = \"130\u00f78=16, 2\";   New = \"106\u00f72=53, 0\" },\n    @{ Old = \"338\u00f77=48, 2\";   New = \"785\u00f74=196, 1\" },\n    @{ Old = \"646\u00f79=71, 7\";   New = \"995\u00f72=497, 1\" },\n    @{ Old = \"488\u00f77=69, 5\";   New = \"113\u00f78=14, 1\" },\n    @{ Old = \"491\u00f73=163, 2\";  New = \"234\u00f78=29, 2\" },\n    @{ Old = \"232\u00f77=33, 1\";   New = \"429\u00f75=85, 4\" },\n    @{ Old = \"215\u00f79=23, 8\";   New = \"224\u00f75=44, 4\" },\n    @{ Old = \"256\u00f74=64, 0\";   New = \"266\u00f77=38, 0\" },\n    @{ Old = \"376\u00f73=125, 1\";  New = \"851\u00f79=94, 5\" },\n    @{ Old = \"550\u00f75=110, 0\";  New = \"556\u00f77=79, 3\" },\n    @{ Old = \"519\u00f74=129, 3\";  New = \"154\u00f74=38, 2\" },\n    @{ Old = \"235\u00f78=29, 3\";   New = \"181\u00f72=90, 1\" },\n    @{ Old = \"537\u00f76=89, 3\";   New = \"811\u00f74=202, 3\" },\n    @{ Old = \"838\u00f75=167, 3\";  New = \"197\u00f76=32, 5\" },\n    @{ Old = \"723\u00f75=144, 3\";  New = \"226\u00f72=113, 0\" },\n    @{ Old = \"675\u00f79=75, 0\";   New = \"421\u00f75=84, 1\" },\n    @{ Old = \"122\u00f78=15, 2\";   New = \"378\u00f75=75, 3\" },\n    @{ Old = \"535\u00f79=59, 4\";   New = \"118\u00f77=16, 6\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll)\n}\n"}
